$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.835.31'
$ws.Cells.Item(2, 5).Value = '  -0.22%  '
$ws.Cells.Item(3, 4).Value = '2.077.75'
$ws.Cells.Item(3, 5).Value = '  -1.21%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).Value = '233.52'
$ws.Cells.Item(5, 5).Value = '  +0.02%  '
$ws.Cells.Item(6, 4).Value = '0.626'
$ws.Cells.Item(6, 5).Value = '  +0.34%  '
$ws.Cells.Item(7, 4).Value = '59.32'
$ws.Cells.Item(7, 5).Value = '  +2.01%  '
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 4).Value = '0.393'
$ws.Cells.Item(10, 4).Value = '0.0791'
$ws.Cells.Item(11, 5).Value = '  +1.58%  '
$ws.Cells.Item(12, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(12, 4).Value = '2.383.92'
$ws.Cells.Item(12, 5).Value = '  -0.87%  '
$ws.Cells.Item(13, 2).Value = 'Chainlink'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(13, 4).Value = '14.80'
$ws.Cells.Item(13, 5).Value = '  +1.55%  '
$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).Value = '21.22'
$ws.Cells.Item(14, 5).Value = '  -0.26%  '
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).Value = '0.774'
$ws.Cells.Item(15, 5).Value = '  +0.18%  '
$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(16, 4).Value = '5.36'
$ws.Cells.Item(16, 5).Value = '  +1.88%  '
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '2.110.21'
$ws.Cells.Item(17, 5).Value = '  +0.52%  '
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '37.732.59'
$ws.Cells.Item(18, 5).Value = '  -0.22%  '
$ws.Cells.Item(19, 2).Value = 'Uniswap'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(19, 4).Value = '6.16'
$ws.Cells.Item(19, 5).Value = '  -0.37%  '
$ws.Cells.Item(20, 2).Value = 'Litecoin'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(20, 4).Value = '71.60'
$ws.Cells.Item(20, 5).Value = '  +0.90%  '
$ws.Cells.Item(21, 2).Value = 'ShibaInu'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(21, 4).Value = '0.0₃0855'
$ws.Cells.Item(21, 5).Value = '  +3.55%  '
$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Value = '228.37'
$ws.Cells.Item(22, 5).Value = '  +0.10%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '
$ws.Cells.Item(24, 2).Value = 'PancakeSwap'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(24, 4).Value = '2.42'
$ws.Cells.Item(24, 5).Value = '  +1.01%  '
$ws.Cells.Item(25, 2).Value = 'Toncoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(25, 4).Value = '2.37'
$ws.Cells.Item(25, 5).Value = '  -1.44%  '
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26, 4).Value = '170.58'
$ws.Cells.Item(26, 5).Value = '  +1.38%  '
$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(27, 4).Value = '9.20'
$ws.Cells.Item(27, 5).Value = '  +2.52%  '
$ws.Cells.Item(28, 2).Value = 'Kaspa'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(28, 4).Value = '0.134'
$ws.Cells.Item(28, 5).Value = '  -4.33%  '
$ws.Cells.Item(29, 2).Value = 'ImmutableX'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(29, 4).Value = '1.42'
$ws.Cells.Item(29, 5).Value = '  -0.56%  '
$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).Value = '19.49'
$ws.Cells.Item(30, 5).Value = '  -0.16%  '
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).Value = '0.121'
$ws.Cells.Item(31, 5).Value = '  +1.36%  '
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(32, 4).Value = '4.73'
$ws.Cells.Item(32, 5).Value = '  +1.61%  '
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 4).Value = '4.75'
$ws.Cells.Item(33, 5).Value = '  +2.89%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '0.0633'
$ws.Cells.Item(34, 5).Value = '  +0.53%  '
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(35, 4).Value = '2.49'
$ws.Cells.Item(35, 5).Value = '  -1.67%  '
$ws.Cells.Item(36, 5).Value = '  -0.94%  '
$ws.Cells.Item(37, 2).Value = 'WEMIXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(37, 4).Value = '1.82'
$ws.Cells.Item(37, 5).Value = '  -0.36%  '
$ws.Cells.Item(38, 2).Value = 'BinanceUSD'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(38, 4).Value = '1.00'
$ws.Cells.Item(38, 5).Value = '  +0.20%  '
$ws.Cells.Item(39, 2).Value = 'THORChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(39, 4).Value = '5.41'
$ws.Cells.Item(39, 5).Value = '  -0.27%  '
$ws.Cells.Item(40, 2).Value = 'Cronos'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(40, 4).Value = '0.0984'
$ws.Cells.Item(40, 5).Value = '  -0.87%  '
$ws.Cells.Item(41, 2).Value = 'Aave'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(41, 4).Value = '99.23'
$ws.Cells.Item(41, 5).Value = '  +1.61%  '
$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(42, 4).Value = '0.0217'
$ws.Cells.Item(42, 5).Value = '  +0.65%  '
$ws.Cells.Item(43, 2).Value = 'HuobiToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(43, 4).Value = '2.88'
$ws.Cells.Item(43, 5).Value = '  -2.22%  '
$ws.Cells.Item(44, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(44, 4).Value = '16.63'
$ws.Cells.Item(44, 5).Value = '  +5.41%  '
$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(45, 4).Value = '1.443.66'
$ws.Cells.Item(45, 5).Value = '  -1.04%  '
$ws.Cells.Item(46, 4).Value = '4.26'
$ws.Cells.Item(46, 5).Value = '  +4.85%  '
$ws.Cells.Item(47, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(47, 4).Value = '1.16'
$ws.Cells.Item(47, 5).Value = '  -0.98%  '
$ws.Cells.Item(48, 2).Value = 'ARBITRUM'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(48, 4).Value = '1.06'
$ws.Cells.Item(48, 5).Value = '  +0.08%  '
$ws.Cells.Item(49, 2).Value = 'FraxShare'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(49, 4).Value = '7.41'
$ws.Cells.Item(49, 5).Value = '  +0.42%  '
$ws.Cells.Item(50, 2).Value = 'MXToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(50, 4).Value = '3.01'
$ws.Cells.Item(50, 5).Value = '  -0.36%  '
$ws.Cells.Item(51, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(51, 4).Value = '2.268.59'
$ws.Cells.Item(51, 5).Value = '  -1.20%  '
